# checkbox auto nuevo y subir factura

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Datos Extraídos"

# Duplicate row 2 (the existing record) down into row 3 before touching row 2,
# so row 3 inherits the same text-typed cell values/format as row 2 (avoids
# Excel's automatic number/date inference when re-typing values like
# "17-JUN-2004", "2020" or "1.248").
$ws.Range("A2:R2").Copy()
$ws.Range("A3:R3").PasteSpecial(-4104)
$ws.Range("A3:R3").Borders.LineStyle = 1

# New row 3 keeps the person's data, but reflects the new vehicle's
# circulation city/department, gender and that it WAS purchased with a bank
# loan (DAVIVIENDA).
$ws.Range("O3").Value = "Risaralda"
$ws.Range("P3").Value = "Pereira"
$ws.Range("Q3").Value = "MASCULINO"
$ws.Range("R3").Value = "DAVIVIENDA"

# Row 2 (the original record) now reflects its own, different, circulation
# city/department + gender, and is NOT bank-financed.
$ws.Range("O2").Value = "Risaralda"
$ws.Range("P2").Value = "Dosquebradas"
$ws.Range("Q2").Value = "MASCULINO"
$ws.Range("R2").Value = "NO"

# Widen column Q slightly (stored width 10 -> 11). The engine's ColumnWidth
# (character units) maps to stored width as stored = ColumnWidth + 5/6, so
# subtract that padding to land exactly on a stored width of 11.
$ws.Columns.Item(17).ColumnWidth = 10.166666666666666
